$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planta Daninha")

# Fill in the new "Resposta" column (I) first: NAO, then SIM (matches the
# shared-string insertion order captured in the target workbook).
$ws.Cells.Item(15, 9).Value = "NÃO"
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 9).Value = "SIM"
}

# Replace the question text in D2 with the new wording.
$ws.Range("D2").Value = " Ocorre escape de plantas daninhas na maioria dos talhões da fazenda?"

# Match the saved selection state.
$ws.Range("D3").Select()
